$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the numeric values in column A (clientNumber) for rows 2-5
$ws.Range("A2:A5").Value = 123456789

# 2. Merge the two rich-text runs of the "See <url>" string in H3 into a single
#    uniformly-colored (blue) run, and give the cell its own style referencing
#    that same blue font, so the whole cell renders consistently.
$ws.Range("H3").Font.Color = 16711680
$ws.Range("H3").Characters(1, 4).Font.Color = 16711680
$ws.Range("H3").Characters(5, 56).Font.Color = 16711680

# 3. Update the hyperlink's display text to include the "See " prefix, in place
#    (iterate so we bind to the existing hyperlink object instead of creating a
#    duplicate one).
foreach ($hl in $ws.Hyperlinks) {
    $hl.TextToDisplay = "See https://en.wikipedia.org/wiki/List_of_Unicode_characters"
}

# 4. Update the active selection to A3:A5 with A3 as the active cell.
$ws.Activate()
$ws.Range("A3:A5").Select()
